$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "time_taken" (column F) timestamps on the existing "data"
#    sheet (rows 2..129) to the re-run values.
# ---------------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")

$newTimestamps = @(
    "2021-10-05 14:33:59.794970",
    "2021-10-05 14:33:59.794979",
    "2021-10-05 14:33:59.794982",
    "2021-10-05 14:33:59.794985",
    "2021-10-05 14:33:59.794988",
    "2021-10-05 14:33:59.794991",
    "2021-10-05 14:33:59.794993",
    "2021-10-05 14:33:59.794996",
    "2021-10-05 14:33:59.794999",
    "2021-10-05 14:33:59.795001",
    "2021-10-05 14:33:59.795004",
    "2021-10-05 14:33:59.795007",
    "2021-10-05 14:33:59.795009",
    "2021-10-05 14:33:59.795012",
    "2021-10-05 14:33:59.795015",
    "2021-10-05 14:33:59.795017",
    "2021-10-05 14:33:59.795020",
    "2021-10-05 14:33:59.795023",
    "2021-10-05 14:33:59.795026",
    "2021-10-05 14:33:59.795028",
    "2021-10-05 14:33:59.795031",
    "2021-10-05 14:33:59.795034",
    "2021-10-05 14:33:59.795036",
    "2021-10-05 14:33:59.795039",
    "2021-10-05 14:33:59.795043",
    "2021-10-05 14:33:59.795045",
    "2021-10-05 14:33:59.795048",
    "2021-10-05 14:33:59.795050",
    "2021-10-05 14:33:59.795053",
    "2021-10-05 14:33:59.795056",
    "2021-10-05 14:33:59.795058",
    "2021-10-05 14:33:59.795061",
    "2021-10-05 14:33:59.795064",
    "2021-10-05 14:33:59.795067",
    "2021-10-05 14:33:59.795070",
    "2021-10-05 14:33:59.795072",
    "2021-10-05 14:33:59.795075",
    "2021-10-05 14:33:59.795077",
    "2021-10-05 14:33:59.795080",
    "2021-10-05 14:33:59.795083",
    "2021-10-05 14:33:59.795086",
    "2021-10-05 14:33:59.795089",
    "2021-10-05 14:33:59.795092",
    "2021-10-05 14:33:59.795094",
    "2021-10-05 14:33:59.795097",
    "2021-10-05 14:33:59.795100",
    "2021-10-05 14:33:59.795103",
    "2021-10-05 14:33:59.795106",
    "2021-10-05 14:33:59.795108",
    "2021-10-05 14:33:59.795111",
    "2021-10-05 14:33:59.795113",
    "2021-10-05 14:33:59.795116",
    "2021-10-05 14:33:59.795119",
    "2021-10-05 14:33:59.795122",
    "2021-10-05 14:33:59.795125",
    "2021-10-05 14:33:59.795127",
    "2021-10-05 14:33:59.795130",
    "2021-10-05 14:33:59.795133",
    "2021-10-05 14:33:59.795136",
    "2021-10-05 14:33:59.795138",
    "2021-10-05 14:33:59.795141",
    "2021-10-05 14:33:59.795144",
    "2021-10-05 14:33:59.795147",
    "2021-10-05 14:33:59.795150",
    "2021-10-05 14:33:59.795153",
    "2021-10-05 14:33:59.795156",
    "2021-10-05 14:33:59.795159",
    "2021-10-05 14:33:59.795162",
    "2021-10-05 14:33:59.795164",
    "2021-10-05 14:33:59.795167",
    "2021-10-05 14:33:59.795170",
    "2021-10-05 14:33:59.795172",
    "2021-10-05 14:33:59.795175",
    "2021-10-05 14:33:59.795178",
    "2021-10-05 14:33:59.795180",
    "2021-10-05 14:33:59.795183",
    "2021-10-05 14:33:59.795188",
    "2021-10-05 14:33:59.795191",
    "2021-10-05 14:33:59.795194",
    "2021-10-05 14:33:59.795197",
    "2021-10-05 14:33:59.795199",
    "2021-10-05 14:33:59.795202",
    "2021-10-05 14:33:59.795205",
    "2021-10-05 14:33:59.795207",
    "2021-10-05 14:33:59.795210",
    "2021-10-05 14:33:59.795213",
    "2021-10-05 14:33:59.795215",
    "2021-10-05 14:33:59.795218",
    "2021-10-05 14:33:59.795220",
    "2021-10-05 14:33:59.795223",
    "2021-10-05 14:33:59.795226",
    "2021-10-05 14:33:59.795228",
    "2021-10-05 14:33:59.795232",
    "2021-10-05 14:33:59.795235",
    "2021-10-05 14:33:59.795238",
    "2021-10-05 14:33:59.795241",
    "2021-10-05 14:33:59.795244",
    "2021-10-05 14:33:59.795246",
    "2021-10-05 14:33:59.795249",
    "2021-10-05 14:33:59.795251",
    "2021-10-05 14:33:59.795254",
    "2021-10-05 14:33:59.795256",
    "2021-10-05 14:33:59.795259",
    "2021-10-05 14:33:59.795262",
    "2021-10-05 14:33:59.795265",
    "2021-10-05 14:33:59.795267",
    "2021-10-05 14:33:59.795270",
    "2021-10-05 14:33:59.795272",
    "2021-10-05 14:33:59.795277",
    "2021-10-05 14:33:59.795280",
    "2021-10-05 14:33:59.795283",
    "2021-10-05 14:33:59.795286",
    "2021-10-05 14:33:59.795289",
    "2021-10-05 14:33:59.795292",
    "2021-10-05 14:33:59.795294",
    "2021-10-05 14:33:59.795297",
    "2021-10-05 14:33:59.795300",
    "2021-10-05 14:33:59.795302",
    "2021-10-05 14:33:59.795305",
    "2021-10-05 14:33:59.795308",
    "2021-10-05 14:33:59.795311",
    "2021-10-05 14:33:59.795313",
    "2021-10-05 14:33:59.795316",
    "2021-10-05 14:33:59.795319",
    "2021-10-05 14:33:59.795321",
    "2021-10-05 14:33:59.795324",
    "2021-10-05 14:33:59.795327",
    "2021-10-05 14:33:59.795330"
)


for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Formula = '="' + $newTimestamps[$i] + '"'
}

# Turn the just-written formulas back into plain literal text values so the
# cells round-trip as inline strings (matching the original authoring tool),
# not formulas.
$usedRange = $dataSheet.Range("F2:F129")
$usedRange.Copy() | Out-Null
$usedRange.PasteSpecial(-4163) | Out-Null  # xlPasteValues

# ---------------------------------------------------------------------------
# 2) Add the new "metadata" worksheet (placed after "data") describing the
#    panel query that produced this export.
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hereditary Spastic Paraplegia - paediatric"
$metaSheet.Range("C2").Value = 317
$metaSheet.Range("D2").Formula = '="1.18"'
$metaSheet.Range("E2").Value = "2021-10-04T05:19:06.692274Z"
$metaSheet.Range("F2").Formula = '="2021-10-05 14:33:59.791653"'
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/317/?format=json"

# Collapse the D2/F2 formulas down to literal text values (same reasoning as
# above: the source data are plain strings, not formulas).
$metaRow2 = $metaSheet.Range("D2:D2")
$metaRow2.Copy() | Out-Null
$metaRow2.PasteSpecial(-4163) | Out-Null

$metaF2 = $metaSheet.Range("F2:F2")
$metaF2.Copy() | Out-Null
$metaF2.PasteSpecial(-4163) | Out-Null

# Apply the same header style (bold + border, centered/top-aligned) already
# used by the "data" sheet's header row / index column, reusing the existing
# style slot instead of minting a new one.
$dataSheet.Range("B1").Copy() | Out-Null
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$dataSheet.Range("A2").Copy() | Out-Null
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

Write-Output "done"
